$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>85 x 98</w:t><w:br/><w:t xml:space="preserve">  9    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>')
$t.Cell(1,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>18 x 87</w:t><w:br/><w:t xml:space="preserve">  8    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p>')
$t.Cell(1,3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>89 x 95</w:t><w:br/><w:t xml:space="preserve">  9    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>')
$t.Cell(2,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>15 x 23</w:t><w:br/><w:t xml:space="preserve">  2    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p>')
$t.Cell(2,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>23 x 58</w:t><w:br/><w:t xml:space="preserve">  5    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>')
$t.Cell(2,3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>60 x 54</w:t><w:br/><w:t xml:space="preserve">  5    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p>')
$t.Cell(3,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>31 x 12</w:t><w:br/><w:t xml:space="preserve">  1    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p>')
$t.Cell(3,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>94 x 25</w:t><w:br/><w:t xml:space="preserve">  2    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p>')
$t.Cell(3,3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>23 x 48</w:t><w:br/><w:t xml:space="preserve">  4    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>')
$t.Cell(4,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>56 x 54</w:t><w:br/><w:t xml:space="preserve">  5    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>')
$t.Cell(4,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>49 x 45</w:t><w:br/><w:t xml:space="preserve">  4    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p>')
$t.Cell(4,3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>73 x 56</w:t><w:br/><w:t xml:space="preserve">  5    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p>')
$t.Cell(5,1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>66 x 36</w:t><w:br/><w:t xml:space="preserve">  3    6</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p>')
$t.Cell(5,2).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>90 x 83</w:t><w:br/><w:t xml:space="preserve">  8    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:br/><w:t>0|    |</w:t></w:r></w:p>')
$t.Cell(5,3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>20 x 19</w:t><w:br/><w:t xml:space="preserve">  1    9</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p>')
